$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Retirement control"
$ws2.Activate()

$ws2.Range("A1").Value = "*This is to control whether retirement and planned plants are used or not"
$ws2.Range("A2").Value = "* moved here so that it sees all scenarios for CAP_BND"
$ws2.Range("B4").Value = "~TFM_INS"
$ws2.Range("B5").Value = "TimeSlice"
$ws2.Range("C5").Value = "LimType"
$ws2.Range("D5").Value = "Attribute"
$ws2.Range("E5").Value = "Year"
$ws2.Range("F5").Value = "Attrib_Cond"
$ws2.Range("G5").Value = "AllRegions"
$ws2.Range("H5").Value = "Pset_PD"
$ws2.Range("I5").Value = "Pset_CI"
$ws2.Range("C6").Value = "LO"
$ws2.Range("D6").Value = "RCAP_BND"
$ws2.Range("E6").Value = 2010
$ws2.Range("F6").Value = "CAP_BND"
$ws2.Range("G6").Value = 0
$ws2.Range("H6").Value = "Existing Elec*,Existing CHP*"
$ws2.Range("C7").Value = "UP"
$ws2.Range("D7").Value = "RCAP_BND"
$ws2.Range("E7").Value = 2030
$ws2.Range("F7").Value = "-CAP_BND"
$ws2.Range("G7").Value = 0
$ws2.Range("H7").Value = "Existing Elec*,Existing CHP*"
$ws2.Range("I7").Value = "ELCSOL,ELCWIN"
$ws2.Range("C8").Value = "UP"
$ws2.Range("D8").Value = "RCAP_BND"
$ws2.Range("E8").Value = 2020
$ws2.Range("F8").Value = "-CAP_BND"
$ws2.Range("G8").Value = 0
$ws2.Range("H8").Value = "Existing Elec*,Existing CHP*"
$ws2.Range("I8").Value = "-ELCSOL,-ELCWIN"
$ws2.Range("C9").Value = "UP"
$ws2.Range("D9").Value = "RCAP_BND"
$ws2.Range("E9").Value = 0
$ws2.Range("F9").Value = "-CAP_BND"
$ws2.Range("G9").Value = 4
$ws2.Range("H9").Value = "Existing Elec*,Existing CHP*"
$ws2.Columns.Item(2).AutoFit() | Out-Null
$ws2.Columns.Item(3).AutoFit() | Out-Null
$ws2.Columns.Item(4).AutoFit() | Out-Null
$ws2.Columns.Item(5).AutoFit() | Out-Null
$ws2.Columns.Item(6).AutoFit() | Out-Null
$ws2.Columns.Item(7).AutoFit() | Out-Null
$ws2.Columns.Item(8).AutoFit() | Out-Null
$ws2.Columns.Item(9).AutoFit() | Out-Null
$ws2.Range("F10").Select()
